$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two data rows (58 -> 56 total data rows incl. header)
$ws.Rows.Item(57).Delete()
$ws.Rows.Item(57).Delete()

# Updated lamda_1 (B), lamda_2 (C), dic_nbre_clients_poisson_2_keys (D)
# and dic_nbre_clients_prob_poisson_2_values (E) for rows 2-56.
# Column A (index) is unchanged.
$data = @(
        @(33.94444444444444,1.95,0,0.124),
        @(33.94444444444444,1.95,2,0.002),
        @(33.94444444444444,1.95,3,0.007),
        @(33.94444444444444,1.95,4,0.01),
        @(33.94444444444444,1.95,5,0.022),
        @(33.94444444444444,1.95,6,0.03),
        @(33.94444444444444,1.95,7,0.048),
        @(33.94444444444444,1.95,8,0.052),
        @(33.94444444444444,1.95,9,0.043),
        @(33.94444444444444,1.95,10,0.041),
        @(33.94444444444444,1.95,11,0.025),
        @(33.94444444444444,1.95,12,0.031),
        @(33.94444444444444,1.95,13,0.022),
        @(33.94444444444444,1.95,14,0.037),
        @(33.94444444444444,1.95,15,0.036),
        @(33.94444444444444,1.95,16,0.035),
        @(33.94444444444444,1.95,17,0.029),
        @(33.94444444444444,1.95,18,0.043),
        @(33.94444444444444,1.95,19,0.023),
        @(33.94444444444444,1.95,20,0.022),
        @(33.94444444444444,1.95,21,0.027),
        @(33.94444444444444,1.95,22,0.019),
        @(33.94444444444444,1.95,23,0.028),
        @(33.94444444444444,1.95,24,0.022),
        @(33.94444444444444,1.95,25,0.019),
        @(33.94444444444444,1.95,26,0.015),
        @(33.94444444444444,1.95,27,0.026),
        @(33.94444444444444,1.95,28,0.018),
        @(33.94444444444444,1.95,29,0.016),
        @(33.94444444444444,1.95,30,0.017),
        @(33.94444444444444,1.95,31,0.011),
        @(33.94444444444444,1.95,32,0.012),
        @(33.94444444444444,1.95,33,0.011),
        @(33.94444444444444,1.95,34,0.007),
        @(33.94444444444444,1.95,35,0.009000000000000001),
        @(33.94444444444444,1.95,36,0.007),
        @(33.94444444444444,1.95,37,0.004),
        @(33.94444444444444,1.95,38,0.008),
        @(33.94444444444444,1.95,39,0.005),
        @(33.94444444444444,1.95,40,0.004),
        @(33.94444444444444,1.95,41,0.004),
        @(33.94444444444444,1.95,42,0.001),
        @(33.94444444444444,1.95,43,0.004),
        @(33.94444444444444,1.95,44,0.004),
        @(33.94444444444444,1.95,45,0.003),
        @(33.94444444444444,1.95,46,0.002),
        @(33.94444444444444,1.95,47,0.001),
        @(33.94444444444444,1.95,48,0.002),
        @(33.94444444444444,1.95,49,0.003),
        @(33.94444444444444,1.95,50,0.002),
        @(33.94444444444444,1.95,51,0.001),
        @(33.94444444444444,1.95,52,0.002),
        @(33.94444444444444,1.95,54,0.001),
        @(33.94444444444444,1.95,55,0.001),
        @(33.94444444444444,1.95,79,0.001)
    )

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
}
